$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'61.849.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = '  +1.27%  '
$ws.Range("D3").Formula = "'3.461.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = '  +2.34%  '
$ws.Range("D4").Formula = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = '  -0.01%  '
$ws.Range("D5").Formula = "'583.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +1.54%  '
$ws.Range("D6").Formula = "'147.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  +7.42%  '
$ws.Range("D7").Formula = "'3.462.37"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = '  +2.41%  '
$ws.Range("E8").Value2 = '  -0.08%  '
$ws.Range("E9").Value2 = '  +1.36%  '
$ws.Range("E10").Value2 = '  +0.15%  '
$ws.Range("E11").Value2 = '  +3.48%  '
$ws.Range("E12").Value2 = '  +2.75%  '
$ws.Range("D13").Formula = "'4.052.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = '  +2.39%  '
$ws.Range("D14").Formula = "'27.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = '  +9.26%  '
$ws.Range("E15").Value2 = '  -0.93%  '
$ws.Range("E16").Value2 = '  +1.59%  '
$ws.Range("D17").Formula = "'3.465.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = '  +2.57%  '
$ws.Range("D18").Formula = "'61.942.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = '  +1.16%  '
$ws.Range("E19").Value2 = '  +8.73%  '
$ws.Range("D20").Formula = "'14.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Formula = "'9.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  +2.90%  '
$ws.Range("D22").Formula = "'389.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  +3.28%  '
$ws.Range("E23").Value2 = '  +2.81%  '
$ws.Range("D24").Formula = "'73.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  +3.87%  '
$ws.Range("E25").Value2 = '  +0.25%  '
$ws.Range("D26").Formula = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = '  -0.22%  '
$ws.Range("E27").Value2 = '  -1.82%  '
$ws.Range("D28").Formula = "'3.600.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = '  +2.27%  '
$ws.Range("E29").Value2 = '  -0.47%  '
$ws.Range("D30").Formula = "'7.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  +4.44%  '
$ws.Range("D31").Formula = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = '  -0.05%  '
$ws.Range("B32").Value2 = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Formula = "'8.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  +1.82%  '
$ws.Range("D33").Formula = "'1.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = '  -9.95%  '
$ws.Range("B34").Value2 = 'PancakeSwap'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Formula = "'2.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  +2.26%  '
$ws.Range("E35").Value2 = '  +0.00%  '
$ws.Range("D36").Formula = "'24.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  +3.47%  '
$ws.Range("D37").Formula = "'3.487.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = '  +2.40%  '
$ws.Range("B39").Value2 = 'NEARProtocol'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Formula = "'5.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  +0.57%  '
$ws.Range("B40").Value2 = 'ImmutableX'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Formula = "'1.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  +1.96%  '
$ws.Range("D41").Formula = "'167.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  +1.45%  '
$ws.Range("D42").Formula = "'0.0785"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  +3.64%  '
$ws.Range("D43").Formula = "'27.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = '  +7.22%  '
$ws.Range("D44").Formula = "'0.808"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  +4.42%  '
$ws.Range("E45").Value2 = '  +3.87%  '
$ws.Range("D46").Formula = "'42.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("B47").Value2 = 'FirstDigitalUSD'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Formula = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  +0.07%  '
$ws.Range("B48").Value2 = 'Stacks'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Formula = "'1.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = '  +1.90%  '
$ws.Range("E49").Value2 = '  -1.57%  '
$ws.Range("D50").Formula = "'2.570.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  +1.56%  '
$ws.Range("E51").Value2 = '  +2.49%  '
